# Refresh cached market-board derived figures (currentAveragePrice*, Leve
# price/profit columns) for the leve rows whose underlying item prices moved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 1751772.4
$ws.Range("I64").Value = 3208652.2
$ws.Range("J64").Value = 3516.65
$ws.Range("K64").Value = 3208652.2
$ws.Range("L64").Value = 3516.65
$ws.Range("M64").Value = -3208404.2
$ws.Range("N64").Value = -4012.65

$ws.Range("H67").Value = 1751772.4
$ws.Range("I67").Value = 3208652.2
$ws.Range("J67").Value = 3516.65
$ws.Range("K67").Value = 3208652.2
$ws.Range("L67").Value = 3516.65
$ws.Range("M67").Value = -3207794.2
$ws.Range("N67").Value = -5232.65

$ws.Range("H76").Value = 3748.4333
$ws.Range("I76").Value = 3738.7144
$ws.Range("J76").Value = 3771.111
$ws.Range("K76").Value = 3738.7144
$ws.Range("L76").Value = 3771.111
$ws.Range("M76").Value = -3423.7144
$ws.Range("N76").Value = -4401.111

$ws.Range("H79").Value = 3748.4333
$ws.Range("I79").Value = 3738.7144
$ws.Range("J79").Value = 3771.111
$ws.Range("K79").Value = 3738.7144
$ws.Range("L79").Value = 3771.111
$ws.Range("M79").Value = -2646.7144
$ws.Range("N79").Value = -5955.111

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11581.512
$ws.Range("I32").Value = 10598.831
$ws.Range("K32").Value = 10598.831
$ws.Range("M32").Value = -10311.831

$ws.Range("H45").Value = 1468059
$ws.Range("I45").Value = 2526986
$ws.Range("J45").Value = 1852.4615
$ws.Range("K45").Value = 2526986
$ws.Range("L45").Value = 1852.4615
$ws.Range("M45").Value = -2526609
$ws.Range("N45").Value = -2606.4615

$ws.Range("H63").Value = 125002200
$ws.Range("I63").Value = 166668770
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 166668770
$ws.Range("L63").Value = 2500
$ws.Range("M63").Value = -166668084
$ws.Range("N63").Value = -3872

$ws.Range("H66").Value = 125002200
$ws.Range("I66").Value = 166668770
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 833343850
$ws.Range("L66").Value = 12500
$ws.Range("M66").Value = -833340418
$ws.Range("N66").Value = -19364

$ws.Range("H88").Value = 2001.75
$ws.Range("I88").Value = 1600
$ws.Range("J88").Value = 2403.5
$ws.Range("K88").Value = 1600
$ws.Range("L88").Value = 2403.5
$ws.Range("M88").Value = -1194
$ws.Range("N88").Value = -3215.5

$ws.Range("H91").Value = 2001.75
$ws.Range("I91").Value = 1600
$ws.Range("J91").Value = 2403.5
$ws.Range("K91").Value = 1600
$ws.Range("L91").Value = 2403.5
$ws.Range("M91").Value = -196
$ws.Range("N91").Value = -5211.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1140
$ws.Range("I99").Value = 1166.6666
$ws.Range("J99").Value = 1100
$ws.Range("K99").Value = 1166.6666
$ws.Range("L99").Value = 1100
$ws.Range("M99").Value = 331.3334
$ws.Range("N99").Value = -4096

$ws.Range("H105").Value = 4035.805
$ws.Range("I105").Value = 2751.818
$ws.Range("J105").Value = 4506.6
$ws.Range("K105").Value = 2751.818
$ws.Range("L105").Value = 4506.6
$ws.Range("M105").Value = -1004.818
$ws.Range("N105").Value = -8000.6

$ws.Range("H107").Value = 1939.3636
$ws.Range("I107").Value = 1666.6666
$ws.Range("J107").Value = 3166.5
$ws.Range("K107").Value = 1666.6666
$ws.Range("L107").Value = 3166.5
$ws.Range("M107").Value = 253.3334
$ws.Range("N107").Value = -7006.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5749301.5
$ws.Range("I31").Value = 2975.7144
$ws.Range("J31").Value = 7577678
$ws.Range("K31").Value = 2975.7144
$ws.Range("L31").Value = 7577678
$ws.Range("M31").Value = -2680.7144
$ws.Range("N31").Value = -7578268

$ws.Range("H34").Value = 5749301.5
$ws.Range("I34").Value = 2975.7144
$ws.Range("J34").Value = 7577678
$ws.Range("K34").Value = 2975.7144
$ws.Range("L34").Value = 7577678
$ws.Range("M34").Value = -2773.7144
$ws.Range("N34").Value = -7578082

$ws.Range("H59").Value = 16098
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 16098
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 16098
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -18388

$ws.Range("H62").Value = 2630
$ws.Range("J62").Value = 2916.6667
$ws.Range("L62").Value = 2916.6667
$ws.Range("N62").Value = -4164.6667

$ws.Range("H65").Value = 2630
$ws.Range("J65").Value = 2916.6667
$ws.Range("L65").Value = 14583.3335
$ws.Range("N65").Value = -20823.3335

$ws.Range("H107").Value = 528.7241
$ws.Range("I107").Value = 361.26666
$ws.Range("J107").Value = 708.1429000000001
$ws.Range("K107").Value = 361.26666
$ws.Range("L107").Value = 708.1429000000001
$ws.Range("M107").Value = 1558.73334
$ws.Range("N107").Value = -4548.1429

$ws.Range("H135").Value = 51950
$ws.Range("J135").Value = 51950
$ws.Range("L135").Value = 51950
$ws.Range("N135").Value = -62090

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 1715.56
$ws.Range("J118").Value = 1857.2727
$ws.Range("L118").Value = 5571.8181
$ws.Range("N118").Value = -8057.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11517.692
$ws.Range("I70").Value = 14860.223
$ws.Range("J70").Value = 3997
$ws.Range("K70").Value = 14860.223
$ws.Range("L70").Value = 3997
$ws.Range("M70").Value = -14590.223
$ws.Range("N70").Value = -4537

$ws.Range("H73").Value = 11517.692
$ws.Range("I73").Value = 14860.223
$ws.Range("J73").Value = 3997
$ws.Range("K73").Value = 14860.223
$ws.Range("L73").Value = 3997
$ws.Range("M73").Value = -13924.223
$ws.Range("N73").Value = -5869

$ws.Range("H80").Value = 12348513
$ws.Range("I80").Value = 37039384
$ws.Range("J80").Value = 3077.5
$ws.Range("K80").Value = 37039384
$ws.Range("L80").Value = 3077.5
$ws.Range("M80").Value = -37038386
$ws.Range("N80").Value = -5073.5

$ws.Range("H83").Value = 12348513
$ws.Range("I83").Value = 37039384
$ws.Range("J83").Value = 3077.5
$ws.Range("K83").Value = 185196920
$ws.Range("L83").Value = 15387.5
$ws.Range("M83").Value = -185191928
$ws.Range("N83").Value = -25371.5

$ws.Range("H133").Value = 70213.766
$ws.Range("J133").Value = 70213.766
$ws.Range("L133").Value = 70213.766
$ws.Range("N133").Value = -80333.766

$ws.Range("H139").Value = 71663
$ws.Range("J139").Value = 71663
$ws.Range("L139").Value = 71663
$ws.Range("N139").Value = -81943
